$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 46 (formatting copied from the row above by default)
$ws.Rows.Item(46).Insert()

# 2. Populate the new row 46 with the MELASONOZ SYRUP product data
$ws.Cells.Item(46, 3).Value = "MELASONOZ   SYRUP"
$ws.Cells.Item(46, 8).Value = "0:0"
$ws.Cells.Item(46, 12).Value = "0"
$ws.Cells.Item(46, 14).Value = "60.00"
$ws.Cells.Item(46, 16).Value = "60.0000"
$ws.Cells.Item(46, 17).Value = "1:0"

# 3. Re-merge the cells for the newly inserted row 46 to match the other product rows
$ws.Range("A46:B46").Merge() | Out-Null
$ws.Range("C46:G46").Merge() | Out-Null
$ws.Range("H46:K46").Merge() | Out-Null
$ws.Range("L46:M46").Merge() | Out-Null
$ws.Range("N46:O46").Merge() | Out-Null

# 4. Fix the row height for the new row (matches the other product rows)
$ws.Rows.Item(46).RowHeight = 25.5

# 5. Restore column A ("م" sequence number) for rows 46..73 to the plain sequential
#    numbering (this column is independent of the product-data shift).
for ($r = 46; $r -le 73; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# 6. Update the grand-total (row 74, was row 73) to add the new product's price
$ws.Cells.Item(74, 16).Value = 4249.3

# 7. Update the printed timestamp in the footer (row 75, was row 74)
$ws.Cells.Item(75, 1).Value = "Saturday, 6 September, 2025 8:52 PM"

# 8. Fix up row heights that do not follow the simple "shift down" behaviour
$ws.Rows.Item(73).RowHeight = 24.75
$ws.Rows.Item(74).RowHeight = 25.5
$ws.Rows.Item(75).RowHeight = 16.5
